$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "hidden" helper date/time columns H and I for rows 6 and 8:
# shift the stored date back by one day (from 2-Nov-2022 to 1-Nov-2022),
# matching rows 5 and 7.
$ws.Range("H6").Value = 44866.333333333336
$ws.Range("I6").Value = 44866.708333333336
$ws.Range("H8").Value = 44866.333333333336
$ws.Range("I8").Value = 44866.708333333336

# Update the active cell selection shown in the sheet view.
$ws.Range("G32").Select()
